$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '63.852.95'
Set-TextValue 2 5 '  +0.06%  '
Set-TextValue 3 4 '2.627.69'
Set-TextValue 3 5 '  +0.16%  '
Set-TextValue 4 5 '  -0.02%  '
Set-TextValue 5 4 '596.14'
Set-TextValue 5 5 '  +0.02%  '
Set-TextValue 6 4 '151.53'
Set-TextValue 6 5 '  +1.14%  '
Set-TextValue 7 5 '  -0.02%  '
Set-TextValue 8 4 '0.588'
Set-TextValue 8 5 '  +0.17%  '
Set-TextValue 9 4 '0.113'
Set-TextValue 9 5 '  +3.91%  '
Set-TextValue 10 5 '  +3.24%  '
Set-TextValue 11 4 '0.395'
Set-TextValue 11 5 '  +3.73%  '
Set-TextValue 12 5 '  +1.01%  '
Set-TextValue 13 4 '27.99'
Set-TextValue 14 4 '3.097.86'
Set-TextValue 14 5 '  +0.16%  '
Set-TextValue 15 4 '63.694.48'
Set-TextValue 15 5 '  +0.07%  '
Set-TextValue 16 4 '0.0000164'
Set-TextValue 16 5 '  +10.44%  '
Set-TextValue 17 4 '2.611.43'
Set-TextValue 17 5 '  +0.33%  '
Set-TextValue 18 4 '12.26'
Set-TextValue 18 5 '  +0.74%  '
Set-TextValue 19 5 '  +4.31%  '
Set-TextValue 20 4 '348.65'
Set-TextValue 20 5 '  -0.07%  '
Set-TextValue 21 5 '  +2.05%  '
Set-TextValue 22 5 '  +0.41%  '
Set-TextValue 23 4 '67.51'
Set-TextValue 23 5 '  +2.05%  '
Set-TextValue 24 4 '1.70'
Set-TextValue 24 5 '  -2.38%  '
Set-TextValue 25 4 '9.23'
Set-TextValue 25 5 '  +0.79%  '
Set-TextValue 26 5 '  +0.12%  '
Set-TextValue 27 5 '  +3.65%  '
Set-TextValue 28 4 '552.14'
Set-TextValue 28 5 '  +1.94%  '
Set-TextValue 29 4 '0.162'
Set-TextValue 29 5 '  -0.76%  '
Set-TextValue 30 5 '  -0.08%  '
Set-TextValue 31 5 '  +1.21%  '
Set-TextValue 32 4 '0.0₃0894'
Set-TextValue 32 5 '  +6.09%  '
Set-TextValue 33 5 '  +3.52%  '
Set-TextValue 34 5 '  +3.99%  '
Set-TextValue 35 4 '6.16'
Set-TextValue 35 5 '  +1.94%  '
Set-TextValue 36 4 '164.68'
Set-TextValue 37 4 '0.418'
Set-TextValue 37 5 '  +2.75%  '
Set-TextValue 38 4 '1.99'
Set-TextValue 39 4 '19.86'
Set-TextValue 39 5 '  +2.64%  '
Set-TextValue 40 5 '  +0.01%  '
Set-TextValue 42 4 '168.19'
Set-TextValue 42 5 '  -0.73%  '
Set-TextValue 43 5 '  +4.55%  '
Set-TextValue 44 4 '23.73'
Set-TextValue 44 5 '  +11.22%  '
Set-TextValue 45 5 '  -1.04%  '
Set-TextValue 46 5 '  +10.59%  '
Set-TextValue 47 4 '0.638'
Set-TextValue 47 5 '  +1.68%  '
Set-TextValue 48 4 '0.0253'
Set-TextValue 48 5 '  +3.51%  '
Set-TextValue 49 4 '0.0970'
Set-TextValue 49 5 '  +0.58%  '
Set-TextValue 50 4 '19.29'
Set-TextValue 50 5 '  +0.92%  '
Set-TextValue 51 4 '0.0₆0231'
Set-TextValue 51 5 '  +17.84%  '
